# Scheduled market-data refresh: update cached Universalis price columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the Leve-profit
# sheets. Values below mirror the latest pull from the data source.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2568.375
$ws.Range("J62").Value = 2874.75
$ws.Range("L62").Value = 2874.75
$ws.Range("N62").Value = -4122.75

$ws.Range("H65").Value = 2568.375
$ws.Range("J65").Value = 2874.75
$ws.Range("L65").Value = 14373.75
$ws.Range("N65").Value = -20613.75

$ws.Range("H98").Value = 1682.3462
$ws.Range("I98").Value = 1489.2084
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 1489.2084
$ws.Range("L98").Value = 4000
$ws.Range("M98").Value = 8.791600000000017
$ws.Range("N98").Value = -6996

$ws.Range("H113").Value = 2601
$ws.Range("J113").Value = 2750
$ws.Range("L113").Value = 2750
$ws.Range("N113").Value = -9258

$ws.Range("H122").Value = 1682.3462
$ws.Range("I122").Value = 1489.2084
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 4467.6252
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2017.6252
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 20000
$ws.Range("J109").Value = 20000
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22774

$ws.Range("H132").Value = 3625.4263
$ws.Range("I132").Value = 3803.2703
$ws.Range("J132").Value = 3351.25
$ws.Range("K132").Value = 11409.8109
$ws.Range("L132").Value = 10053.75
$ws.Range("M132").Value = -8879.8109
$ws.Range("N132").Value = -15113.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3315.0908
$ws.Range("I86").Value = 2683.25
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 2683.25
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1560.25
$ws.Range("N86").Value = -7246

$ws.Range("H89").Value = 3315.0908
$ws.Range("I89").Value = 2683.25
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 13416.25
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -7800.25
$ws.Range("N89").Value = -36232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 60.25
$ws.Range("I7").Value = 53.76923
$ws.Range("J7").Value = 72.28570999999999
$ws.Range("K7").Value = 53.76923
$ws.Range("L7").Value = 72.28570999999999
$ws.Range("M7").Value = 59.23077
$ws.Range("N7").Value = -298.28571

$ws.Range("H16").Value = 1942.8572
$ws.Range("I16").Value = 1120
$ws.Range("K16").Value = 1120
$ws.Range("M16").Value = -833

$ws.Range("H22").Value = 334.46155
$ws.Range("I22").Value = 320.66666
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 320.66666
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 29.33334000000002
$ws.Range("N22").Value = -1200

$ws.Range("H50").Value = 16066
$ws.Range("J50").Value = 16066
$ws.Range("L50").Value = 16066
$ws.Range("N50").Value = -17316

$ws.Range("H51").Value = 10916
$ws.Range("J51").Value = 11144.25
$ws.Range("L51").Value = 11144.25
$ws.Range("N51").Value = -12616.25

$ws.Range("H59").Value = 29813.1
$ws.Range("J59").Value = 30891.889
$ws.Range("L59").Value = 30891.889
$ws.Range("N59").Value = -33181.889

$ws.Range("H60").Value = 15607.875
$ws.Range("I60").Value = 5700
$ws.Range("J60").Value = 18910.5
$ws.Range("K60").Value = 5700
$ws.Range("L60").Value = 18910.5
$ws.Range("M60").Value = -5189
$ws.Range("N60").Value = -19932.5

$ws.Range("H61").Value = 10916
$ws.Range("J61").Value = 11144.25
$ws.Range("L61").Value = 11144.25
$ws.Range("N61").Value = -11840.25

$ws.Range("H62").Value = 22226208
$ws.Range("I62").Value = 4282.857
$ws.Range("J62").Value = 41670390
$ws.Range("K62").Value = 4282.857
$ws.Range("L62").Value = 41670390
$ws.Range("M62").Value = -3658.857
$ws.Range("N62").Value = -41671638

$ws.Range("H65").Value = 22226208
$ws.Range("I65").Value = 4282.857
$ws.Range("J65").Value = 41670390
$ws.Range("K65").Value = 21414.285
$ws.Range("L65").Value = 208351950
$ws.Range("M65").Value = -18294.285
$ws.Range("N65").Value = -208358190

$ws.Range("H68").Value = 28748.334
$ws.Range("J68").Value = 28748.334
$ws.Range("L68").Value = 28748.334
$ws.Range("N68").Value = -30246.334

$ws.Range("H71").Value = 28748.334
$ws.Range("J71").Value = 28748.334
$ws.Range("L71").Value = 86245.00199999999
$ws.Range("N71").Value = -93733.00199999999

$ws.Range("H113").Value = 1942.8572
$ws.Range("I113").Value = 1120
$ws.Range("K113").Value = 1120
$ws.Range("M113").Value = 1050

$ws.Range("H122").Value = 1034.2142
$ws.Range("I122").Value = 617.9
$ws.Range("J122").Value = 2075
$ws.Range("K122").Value = 1853.7
$ws.Range("L122").Value = 6225
$ws.Range("M122").Value = 596.3000000000002
$ws.Range("N122").Value = -11125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 874.86206
$ws.Range("I122").Value = 499.42105
$ws.Range("J122").Value = 1588.2
$ws.Range("K122").Value = 4494.78945
$ws.Range("L122").Value = 14293.8
$ws.Range("M122").Value = -2044.78945
$ws.Range("N122").Value = -19193.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8648655
$ws.Range("I70").Value = 11337534
$ws.Range("J70").Value = 5828.5713
$ws.Range("K70").Value = 11337534
$ws.Range("L70").Value = 5828.5713
$ws.Range("M70").Value = -11337264
$ws.Range("N70").Value = -6368.5713

$ws.Range("H73").Value = 8648655
$ws.Range("I73").Value = 11337534
$ws.Range("J73").Value = 5828.5713
$ws.Range("K73").Value = 11337534
$ws.Range("L73").Value = 5828.5713
$ws.Range("M73").Value = -11336598
$ws.Range("N73").Value = -7700.5713

$ws.Range("H132").Value = 58151.223
$ws.Range("I132").Value = 68780.664
$ws.Range("J132").Value = 5004
$ws.Range("K132").Value = 206341.992
$ws.Range("L132").Value = 15012
$ws.Range("M132").Value = -203811.992
$ws.Range("N132").Value = -20072

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 724.75
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 849.5
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 849.5
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -1439.5

$ws.Range("H27").Value = 724.75
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 849.5
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 849.5
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -1063.5

$ws.Range("H46").Value = 2038.2
$ws.Range("I46").Value = 1736.4546
$ws.Range("J46").Value = 2275.2856
$ws.Range("K46").Value = 1736.4546
$ws.Range("L46").Value = 2275.2856
$ws.Range("M46").Value = -1548.4546
$ws.Range("N46").Value = -2651.2856

$ws.Range("H55").Value = 290.94446
$ws.Range("I55").Value = 256.92856
$ws.Range("J55").Value = 410
$ws.Range("K55").Value = 256.92856
$ws.Range("L55").Value = 410
$ws.Range("M55").Value = -83.92856
$ws.Range("N55").Value = -756

$ws.Range("H61").Value = 37039684
$ws.Range("I61").Value = 2902
$ws.Range("K61").Value = 2902
$ws.Range("M61").Value = -2700

$ws.Range("H68").Value = 1262.8823
$ws.Range("I68").Value = 1315.9333
$ws.Range("J68").Value = 865
$ws.Range("K68").Value = 1315.9333
$ws.Range("L68").Value = 865
$ws.Range("M68").Value = -566.9332999999999
$ws.Range("N68").Value = -2363

$ws.Range("H71").Value = 1262.8823
$ws.Range("I71").Value = 1315.9333
$ws.Range("J71").Value = 865
$ws.Range("K71").Value = 6579.666499999999
$ws.Range("L71").Value = 4325
$ws.Range("M71").Value = -2835.666499999999
$ws.Range("N71").Value = -11813

$ws.Range("H113").Value = 37039684
$ws.Range("I113").Value = 2902
$ws.Range("K113").Value = 2902
$ws.Range("M113").Value = -732

$ws.Range("H132").Value = 7651.121
$ws.Range("I132").Value = 8185.185
$ws.Range("J132").Value = 5247.8335
$ws.Range("K132").Value = 24555.555
$ws.Range("L132").Value = 15743.5005
$ws.Range("M132").Value = -22025.555
$ws.Range("N132").Value = -20803.5005

